$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Notes text for the "Combat Armor" row (row 3)
$ws.Range("G3").Value = "Resiliant L2: Soaks one wound level from ballistic or energy damage, unless weapon AP +2 or more."

# Update Strength column (B) values
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 0

# Update the active cell selection
$ws.Range("I23").Select()
